$d = $word.ActiveDocument
$d.TrackRevisions = $false

# Locate the target paragraph (keystore guidance paragraph under "Important Notes")
$para = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Release 2.22.0*updated the*keystore*") {
        $para = $p
    }
}

# --- Step 1: "Release 2.22.0" -> "T" (first run) ---
$rng = $d.Range($para.Range.Start, $para.Range.End)
$rng.Find.Execute("Release 2.22.0", $true, $false, $false, $false, $false, $true, 1, $false, "T", 1)

# --- Step 2: " updated the " -> "he " (second run) ---
$rng = $d.Range($para.Range.Start, $para.Range.End)
$rng.Find.Execute(" updated the ", $true, $false, $false, $false, $false, $true, 1, $false, "he ", 1)

# --- Step 3: insert new run " has been updated in this Release" right after "keystore" ---
$rng = $d.Range($para.Range.Start, $para.Range.End)
$rng.Find.Execute("keystore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $d.Range($rng.End, $rng.End)
$insertPoint.InsertAfter(" has been updated in this Release")

# --- Step 4: ", " -> " " (run right after "command line utilities (CLU)") ---
$rng = $d.Range($para.Range.Start, $para.Range.End)
$rng.Find.Execute("(CLU), ", $true, $false, $false, $false, $false, $true, 1, $false, "(CLU) ", 1)

# --- Step 5: remove "before running any commands, " entirely ---
$rng = $d.Range($para.Range.Start, $para.Range.End)
$rng.Find.Execute("before running any commands, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1)

# --- Step 6: " from GitHub master " -> " from GitHub master" ---
$rng = $d.Range($para.Range.Start, $para.Range.End)
$rng.Find.Execute(" from GitHub master ", $true, $false, $false, $false, $false, $true, 1, $false, " from GitHub master", 1)

# --- Step 7: "if you did not update it after Release 2.22.0" -> " before running any commands" ---
$rng = $d.Range($para.Range.Start, $para.Range.End)
$rng.Find.Execute("if you did not update it after Release 2.22.0", $true, $false, $false, $false, $false, $true, 1, $false, " before running any commands", 1)
